{"js": "// Replace the two placeholder \"Special Control Area\" images with\n// hyperlinked text runs pointing at their real-world image URLs\n// (matching the pStyle=\"BodyText\" paragraphs that previously held the\n// <w:drawing> inline pictures).\nconst body = context.document.body;\n\nconst urls = [\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-1.jpg?h=416&w=750\",\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-2.jpg?h=383&w=750\",\n];\n\nfor (const url of urls) {\n  // Re-query each time: once the first picture's range is replaced with\n  // text, the remaining inline-picture collection shifts down by one.\n  const pics = body.inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n\n  if (pics.items.length === 0) {\n    throw new Error(\"Expected an inline picture to replace with hyperlink: \" + url);\n  }\n\n  const rng = pics.items[0].getRange();\n  // Swap the picture out for its URL as plain text, then turn that text\n  // into a real hyperlink (Word automatically applies the \"Hyperlink\"\n  // character style to the run, just like the authoring UI would).\n  rng.insertText(url, Word.InsertLocation.replace);\n  rng.hyperlink = url;\n  await context.sync();\n}\n", "ps1": "# Replace the two placeholder \"Special Control Area\" images with\n# hyperlinked text runs pointing at their real-world image URLs\n# (matching the pStyle=\"BodyText\" paragraphs that previously held the\n# inline pictures).\n$d = $word.ActiveDocument\n\n$urls = @(\n    \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-1.jpg?h=416&w=750\",\n    \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-2.jpg?h=383&w=750\"\n)\n\nforeach ($url in $urls) {\n    # Re-fetch InlineShapes each time: once a picture is removed/replaced,\n    # the remaining shapes shift down by one.\n    $shapes = $d.InlineShapes\n    if ($shapes.Count -eq 0) {\n        throw \"Expected an inline picture to replace with hyperlink: $url\"\n    }\n\n    $shape = $shapes.Item(1)\n    $rng = $shape.Range\n    $shape.Delete()\n\n    # Add the hyperlink at the now-collapsed range; TextToDisplay makes the\n    # visible run text equal to the URL itself, same as typing/pasting it\n    # and letting AutoFormat turn it into a hyperlink.\n    $d.Hyperlinks.Add($rng, $url, $null, $null, $url) | Out-Null\n}\n"}
